$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a row for FiordlandNZ_ROI1 above the current row 4 (Merauke_ROI1),
# shifting everything from row 4 down by one.
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "FiordlandNZ_ROI1"
$ws.Range("C4").Value = 100

# Insert a row for Patagonia_ROI1 above what is now row 7 (Rwenzori_ROI1),
# shifting everything from row 7 down by one.
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Patagonia_ROI1"
$ws.Range("C7").Value = 500

# Renumber the Index column (A) sequentially for the whole table.
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13

$ws.Range("E10").Clear()

$ws.Range("C7").Select() | Out-Null
